$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F, matching style of existing headers (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Per-row timestamps (rows 2-40) captured at the time of the panel export
$timestamps = @(
    "2021-10-05 10:50:54.974566",
    "2021-10-05 10:50:54.974576",
    "2021-10-05 10:50:54.974579",
    "2021-10-05 10:50:54.974582",
    "2021-10-05 10:50:54.974584",
    "2021-10-05 10:50:54.974587",
    "2021-10-05 10:50:54.974589",
    "2021-10-05 10:50:54.974592",
    "2021-10-05 10:50:54.974595",
    "2021-10-05 10:50:54.974597",
    "2021-10-05 10:50:54.974600",
    "2021-10-05 10:50:54.974602",
    "2021-10-05 10:50:54.974605",
    "2021-10-05 10:50:54.974608",
    "2021-10-05 10:50:54.974610",
    "2021-10-05 10:50:54.974613",
    "2021-10-05 10:50:54.974616",
    "2021-10-05 10:50:54.974619",
    "2021-10-05 10:50:54.974621",
    "2021-10-05 10:50:54.974624",
    "2021-10-05 10:50:54.974627",
    "2021-10-05 10:50:54.974629",
    "2021-10-05 10:50:54.974632",
    "2021-10-05 10:50:54.974634",
    "2021-10-05 10:50:54.974638",
    "2021-10-05 10:50:54.974640",
    "2021-10-05 10:50:54.974643",
    "2021-10-05 10:50:54.974645",
    "2021-10-05 10:50:54.974648",
    "2021-10-05 10:50:54.974651",
    "2021-10-05 10:50:54.974653",
    "2021-10-05 10:50:54.974656",
    "2021-10-05 10:50:54.974659",
    "2021-10-05 10:50:54.974662",
    "2021-10-05 10:50:54.974664",
    "2021-10-05 10:50:54.974667",
    "2021-10-05 10:50:54.974669",
    "2021-10-05 10:50:54.974672",
    "2021-10-05 10:50:54.974674"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
